$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely; everything else shifts left by one column
$ws.Range("A1").EntireColumn.Delete()
